$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format of the last filled data row (28) down into the
# three newly-populated rows (29-31) so date/text formatting matches.
$ws.Range("A28:G28").Copy()
$ws.Range("A29:G29").PasteSpecial(-4122)
$ws.Range("A30:G30").PasteSpecial(-4122)
$ws.Range("A31:G31").PasteSpecial(-4122)

# Row 29 - 2/13/2020
$ws.Range("A29").Value = 43874
$ws.Range("B29").Value = "5:00 - 7:50 pm"
$ws.Range("C29").Value = "Myself"
$ws.Range("D29").Value = "Take the midterm, and learn about stakeholders"
$ws.Range("E29").Value = "Learnt about various stakeholders of the system"
$ws.Range("F29").Value = "Stakeholders of a system are not just limited to people who can affect or be affected by the system, stakeholders may also include organizations and people like donors, maintainers, users etc"
$ws.Range("G29").Value = "Relaxed after the midterm"

# Row 30 - 2/16/2020
$ws.Range("A30").Value = 43877
$ws.Range("B30").Value = "11:00 am - 2:00 pm"
$ws.Range("C30").Value = "Vaishakhi,Anjana"
$ws.Range("D30").Value = "Resubmit homework 1"
$ws.Range("E30").Value = "Resubmitted homework 1"
$ws.Range("F30").Value = "Contrary to our belief that our packet was strong, we realized the importance of ""making it easy for a third person to read"" our code. Slowly but surely we are understanding the need of proper and structured documentation, comments, diagrams that explain our approach. "
$ws.Range("G30").Value = "Neutral"

# Row 31 - 2/19/2020
$ws.Range("A31").Value = 43880
$ws.Range("B31").Value = "11:00 pm - 12:30 am"
$ws.Range("C31").Value = "Vaihakhi,Anjana"
$ws.Range("D31").Value = "Work on high level perspective(the big picture)"
$ws.Range("E31").Value = "Identified stakeholders, system domain, and functionality"
$ws.Range("F31").Value = "Up untill now we have dived into the code. This assignment let us see our system from a broader perspective, By identifying the different kinds of stakeholder, I believe we can understand for instance what kind of organizations would use this system(H2), or what will people benefit out of it. We knew the domain of H2. We talked and understand how h2 is unique from other database options out there. "
$ws.Range("G31").Value = "Relaxed"

# Update the window/selection state to match the edited view.
$ws.Application.ActiveWindow.ScrollRow = 30
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F33").Select()

Write-Output "done"
